# Development Log.xlsx edit
# Commit: "TileMap() Refactor so now prints key value pairs correctly."
#
# A new log entry ("Critical bug fix") is inserted at the top of the data
# table (row 4) on the "Dev Log" sheet. The existing entries (previously in
# rows 4-11) all shift down one row (to rows 5-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev Log")

# --- Shift existing data rows 4-11 down to 5-12 (bottom-up so we never
#     clobber a row before it has been copied). Range.Copy preserves the
#     exact style index of every cell (incl. quotePrefix flags), so no new
#     cellXfs entries get minted in styles.xml.
for ($r = 11; $r -ge 4; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("B" + $srcRow + ":G" + $srcRow)
    $dst = $ws.Range("B" + $dstRow + ":G" + $dstRow)
    $src.Copy($dst)
}

# --- New entry content for row 4 -------------------------------------------
# NOTE: order matters for how new shared strings get appended - write F, then
# D, then E to match the authoring order of the source edit.
$ws.Cells.Item(4, 6).Value = "After stepping away from computer for a bit, I came back to realise i had just created the biggest bug in the entire game. For some reason, i had transposed the keyxalue pairs from the data on the breif. so in stead of 8 x [A1], it was showing 1 x [A8] - 1 tile of A8 worth 8 points ! his was the reason why everything was printing wrong, and i was going down such a rabbit whole not realising this.`nCould not believe this how a critical error, so crucial to the game could be missed. I realise this is where the benefit of working in teams or in pairs helps. Sometimes one's eyes just do not spot the most glaringly obvious of bugs."
$ws.Cells.Item(4, 4).Value = "Critical bug fix"
$ws.Cells.Item(4, 5).Value = "Key-value pairs transposed. Step Back Away from the computer !"

$ws.Cells.Item(4, 3).Value = 0.95208333333333328
$ws.Cells.Item(4, 7).Value = 0.52

# Re-assigning .Value resets any quote-prefix style flag, so copy the
# formatting back in from row 5 (identical styling - it used to be row 4's
# original formatting before the shift above) without touching the values.
$ws.Range("B5:G5").Copy()
$ws.Range("B4:G4").PasteSpecial(-4122) # xlPasteFormats

# --- Conditional formatting ranges now cover one additional row ------------
$cfs = $ws.Cells.FormatConditions
$cf1 = $cfs.Item(1)
$cf1.ModifyAppliesToRange($ws.Range("B4:F12"))
$cf2 = $cfs.Item(2)
$cf2.ModifyAppliesToRange($ws.Range("G4:G12"))

# --- View: active cell moves to F4, top of the visible window moves to A2 --
$ws.Activate()
$ws.Range("F4").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Inserted new Dev Log entry at row 4"
